$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range (A2:D9) so stale shared strings are dropped
$ws.Range("A2:D9").ClearContents()

# Write the refreshed competitor dataset (rows 2-39)
$ws.Cells.Item(2, 1).Value2 = 'Brgy. Poblacion, Bustos, Bulacan (near Mercury Drug)'
$ws.Cells.Item(2, 2).Value2 = 14.8716279
$ws.Cells.Item(2, 3).Value2 = 120.8616286
$ws.Cells.Item(2, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(3, 1).Value2 = 'Sta. Rita, Guiguinto, Bulacan'
$ws.Cells.Item(3, 2).Value2 = 14.8514981
$ws.Cells.Item(3, 3).Value2 = 120.8158673
$ws.Cells.Item(3, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(4, 1).Value2 = 'A. Mabini St., Malolos, Bulacan'
$ws.Cells.Item(4, 2).Value2 = 14.8003883
$ws.Cells.Item(4, 3).Value2 = 120.9237343
$ws.Cells.Item(4, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(5, 1).Value2 = 'Marilao Public Market, Marilao, Bulacan'
$ws.Cells.Item(5, 2).Value2 = 14.8364861
$ws.Cells.Item(5, 3).Value2 = 120.7869787
$ws.Cells.Item(5, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(6, 1).Value2 = 'Poblacion, San Rafael, Bulacan'
$ws.Cells.Item(6, 2).Value2 = 14.9128029
$ws.Cells.Item(6, 3).Value2 = 120.7665534
$ws.Cells.Item(6, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(7, 1).Value2 = 'Santa Maria, Bulacan (exact pin on Waze)'
$ws.Cells.Item(7, 2).Value2 = 14.8859515
$ws.Cells.Item(7, 3).Value2 = 120.8589558
$ws.Cells.Item(7, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(8, 1).Value2 = 'San Vicente, Santa Maria, Bulacan'
$ws.Cells.Item(8, 2).Value2 = 14.8686773
$ws.Cells.Item(8, 3).Value2 = 120.8021336
$ws.Cells.Item(8, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(9, 1).Value2 = 'Sapang Palay II, San Jose del Monte, Bulacan'
$ws.Cells.Item(9, 2).Value2 = 14.8233368
$ws.Cells.Item(9, 3).Value2 = 120.9538739
$ws.Cells.Item(9, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(10, 2).Value2 = 14.7662236
$ws.Cells.Item(10, 3).Value2 = 120.996634
$ws.Cells.Item(10, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(11, 2).Value2 = 14.8402062
$ws.Cells.Item(11, 3).Value2 = 120.7404268
$ws.Cells.Item(11, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(12, 2).Value2 = 14.7443385
$ws.Cells.Item(12, 3).Value2 = 120.9707882
$ws.Cells.Item(12, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(13, 2).Value2 = 14.8280513
$ws.Cells.Item(13, 3).Value2 = 120.8763154
$ws.Cells.Item(13, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(14, 2).Value2 = 14.7812707
$ws.Cells.Item(14, 3).Value2 = 120.9353567
$ws.Cells.Item(14, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(15, 2).Value2 = 14.7624177
$ws.Cells.Item(15, 3).Value2 = 120.9483066
$ws.Cells.Item(15, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(16, 2).Value2 = 14.8534585
$ws.Cells.Item(16, 3).Value2 = 120.8607446
$ws.Cells.Item(16, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(17, 2).Value2 = 14.8407169
$ws.Cells.Item(17, 3).Value2 = 120.861446
$ws.Cells.Item(17, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(18, 2).Value2 = 14.8940707
$ws.Cells.Item(18, 3).Value2 = 120.7790445
$ws.Cells.Item(18, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(19, 2).Value2 = 14.8196486
$ws.Cells.Item(19, 3).Value2 = 120.9042869
$ws.Cells.Item(19, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(20, 2).Value2 = 14.8179281
$ws.Cells.Item(20, 3).Value2 = 120.9059894
$ws.Cells.Item(20, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(21, 2).Value2 = 14.7570098
$ws.Cells.Item(21, 3).Value2 = 120.9541807
$ws.Cells.Item(21, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(22, 2).Value2 = 14.7577127
$ws.Cells.Item(22, 3).Value2 = 120.963266
$ws.Cells.Item(22, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(23, 2).Value2 = 14.7926526
$ws.Cells.Item(23, 3).Value2 = 120.8791063
$ws.Cells.Item(23, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(24, 2).Value2 = 14.7353912
$ws.Cells.Item(24, 3).Value2 = 120.9575609
$ws.Cells.Item(24, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(25, 2).Value2 = 14.9022959
$ws.Cells.Item(25, 3).Value2 = 120.849436
$ws.Cells.Item(25, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(26, 2).Value2 = 14.887753
$ws.Cells.Item(26, 3).Value2 = 120.9669466
$ws.Cells.Item(26, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(27, 2).Value2 = 15.0787116
$ws.Cells.Item(27, 3).Value2 = 120.9406643
$ws.Cells.Item(27, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(28, 2).Value2 = 14.8862936
$ws.Cells.Item(28, 3).Value2 = 120.8675936
$ws.Cells.Item(28, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(29, 2).Value2 = 14.8198267
$ws.Cells.Item(29, 3).Value2 = 120.9630624
$ws.Cells.Item(29, 4).Value2 = 'Chooks to Go'

$ws.Cells.Item(30, 2).Value2 = 14.886931
$ws.Cells.Item(30, 3).Value2 = 120.8679547
$ws.Cells.Item(30, 4).Value2 = 'Baliwag Lechon'

$ws.Cells.Item(31, 2).Value2 = 14.839252
$ws.Cells.Item(31, 3).Value2 = 120.8599118
$ws.Cells.Item(31, 4).Value2 = 'Baliwag Lechon'

$ws.Cells.Item(32, 2).Value2 = 14.7688442
$ws.Cells.Item(32, 3).Value2 = 120.9967441
$ws.Cells.Item(32, 4).Value2 = 'Baliwag Lechon'

$ws.Cells.Item(33, 2).Value2 = 14.7339547
$ws.Cells.Item(33, 3).Value2 = 120.9611707
$ws.Cells.Item(33, 4).Value2 = 'Baliwag Lechon'

$ws.Cells.Item(34, 2).Value2 = 14.7353173
$ws.Cells.Item(34, 3).Value2 = 120.9612501
$ws.Cells.Item(34, 4).Value2 = 'Baliwag Lechon'

$ws.Cells.Item(35, 2).Value2 = 14.7476245
$ws.Cells.Item(35, 3).Value2 = 120.9732312
$ws.Cells.Item(35, 4).Value2 = 'Baliwag Lechon'

$ws.Cells.Item(36, 2).Value2 = 14.8893169
$ws.Cells.Item(36, 3).Value2 = 120.8673315
$ws.Cells.Item(36, 4).Value2 = 'Lechon Manok ni Sr. Pedro'

$ws.Cells.Item(37, 2).Value2 = 14.7661597
$ws.Cells.Item(37, 3).Value2 = 120.9956438
$ws.Cells.Item(37, 4).Value2 = 'Lechon Manok ni Sr. Pedro'

$ws.Cells.Item(38, 2).Value2 = 14.822618
$ws.Cells.Item(38, 3).Value2 = 120.9532644
$ws.Cells.Item(38, 4).Value2 = 'Lechon Manok ni Sr. Pedro'

$ws.Cells.Item(39, 2).Value2 = 14.8658526
$ws.Cells.Item(39, 3).Value2 = 120.8597823
$ws.Cells.Item(39, 4).Value2 = 'Lechon Manok ni Sr. Pedro'

# Restore the active-cell selection recorded in the workbook view
$ws.Range("A8").Select()